# Adiciona a seção "Exceções" (título + 4 itens de lista) entre o item
# "GameMode" e o título "Funcionalidades adicionais".

$d = $word.ActiveDocument

# Localiza o parágrafo "GameMode" (último item da lista anterior a
# "Funcionalidades adicionais") de forma robusta, procurando pelo texto.
$gameModeIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "GameMode*") {
        $gameModeIndex = $i
    }
}

$pGameMode = $d.Paragraphs.Item($gameModeIndex)

# --- Cria os 5 novos parágrafos vazios logo após "GameMode" ---
$pGameMode.Range.InsertParagraphAfter()                          # heading "Exceções"
$pHeading = $d.Paragraphs.Item($gameModeIndex + 1)
$pHeading.Range.InsertParagraphAfter()                           # FullInventoryException
$pExc1 = $d.Paragraphs.Item($gameModeIndex + 2)
$pExc1.Range.InsertParagraphAfter()                              # InsufficientFundsException
$pExc2 = $d.Paragraphs.Item($gameModeIndex + 3)
$pExc2.Range.InsertParagraphAfter()                              # InvalidOperationException
$pExc3 = $d.Paragraphs.Item($gameModeIndex + 4)
$pExc3.Range.InsertParagraphAfter()                              # InvalidTypeException
$pExc4 = $d.Paragraphs.Item($gameModeIndex + 5)

# --- Parágrafo de título "Exceções" ---
$pHeading.Range.ListFormat.RemoveNumbers()
$pHeading.Style = $d.Styles.Item("Normal")
$pHeading.Alignment = 1
$pHeading.Range.Font.Bold = 1
$pHeading.Range.Font.BoldBi = 1
$pHeading.Range.Font.Size = 12
$pHeading.Range.Font.SizeBi = 12
$pHeading.Range.Text = "Exceções"

function Add-ExceptionBullet {
    param(
        $para,
        [string]$name,
        [string]$sep,
        [string]$desc
    )

    $para.Range.Text = $name + $sep + $desc

    $start = $para.Range.Start
    $nameEnd = $start + $name.Length
    $sepEnd = $nameEnd + $sep.Length
    $descEnd = $sepEnd + $desc.Length

    $boldRange = $d.Range($start, $sepEnd)
    $boldRange.Font.Bold = 1
    $boldRange.Font.BoldBi = 1

    $descRange = $d.Range($sepEnd, $descEnd)
    $descRange.Font.Bold = 0
    $descRange.Font.BoldBi = 0
}

Add-ExceptionBullet $pExc1 "FullInventoryException" "-" " Exceção lançada quando o inventário do jogador está cheio e ele tenta adicionar alguma coisa."
Add-ExceptionBullet $pExc2 "InsufficientFundsException" " –" " Exceção lançada quando o jogador não possui saldo suficiente para fazer uma compra."
Add-ExceptionBullet $pExc3 "InvalidOperationException" " –" " Exceção lançada quando o jogador tenta realizar uma operação inválida."
Add-ExceptionBullet $pExc4 "InvalidTypeException" " – " "Exceção lançada quando o usuário tenta criar um objeto com um tipo inváldido."
